# Updates cryptos list values (price/volume columns) to match the latest
# scrape, including a position swap between ShibaInu and Avalanche rows.
# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (e.g. "1.000", "237.80", "23.30", "102.70", "0.000007861") are prefixed
# with a leading apostrophe so they are kept as literal text, matching the
# original inline-string cell contents.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.276.80"
$ws.Range("E2").Value = "  +0.49%  "

$ws.Range("D3").Value = "1.861.00"
$ws.Range("E3").Value = "  +0.63%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "0.7028"

$ws.Range("D6").Value = "'237.80"
$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").Value = "0.08238"
$ws.Range("E8").Value = "  +9.70%  "

$ws.Range("D9").Value = "0.3043"
$ws.Range("E9").Value = "  -0.37%  "

$ws.Range("D10").Value = "'23.30"
$ws.Range("E10").Value = "  -0.43%  "

$ws.Range("D11").Value = "0.08171"
$ws.Range("E11").Value = "  +0.39%  "

$ws.Range("D12").Value = "1.869.31"
$ws.Range("E12").Value = "  +0.59%  "

$ws.Range("D13").Value = "0.7169"
$ws.Range("E13").Value = "  -1.07%  "

$ws.Range("D14").Value = "5.179"
$ws.Range("E14").Value = "  -0.89%  "

$ws.Range("D15").Value = "89.31"
$ws.Range("E15").Value = "  +0.08%  "

$ws.Range("D16").Value = "29.298.81"
$ws.Range("E16").Value = "  +0.50%  "

$ws.Range("D17").Value = "5.785"

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.000007861"
$ws.Range("E18").Value = "  +2.30%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "13.39"
$ws.Range("E19").Value = "  +2.36%  "

$ws.Range("D20").Value = "237.61"
$ws.Range("E20").Value = "  -0.83%  "

$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("D22").Value = "2.108.65"
$ws.Range("E22").Value = "  +0.57%  "

$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("D24").Value = "7.465"
$ws.Range("E24").Value = "  -1.30%  "

$ws.Range("D25").Value = "162.04"
$ws.Range("E25").Value = "  +0.44%  "

$ws.Range("D26").Value = "8.987"
$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("E27").Value = "  -2.08%  "

$ws.Range("D28").Value = "18.12"
$ws.Range("E28").Value = "  +0.48%  "

$ws.Range("E29").Value = "  +1.94%  "

$ws.Range("D30").Value = "1.433"
$ws.Range("E30").Value = "  +3.36%  "

$ws.Range("D31").Value = "4.435"
$ws.Range("E31").Value = "  -2.82%  "

$ws.Range("D32").Value = "1.483"
$ws.Range("E32").Value = "  -0.84%  "

$ws.Range("D33").Value = "4.063"
$ws.Range("E33").Value = "  +1.37%  "

$ws.Range("D34").Value = "0.05224"
$ws.Range("E34").Value = "  +0.83%  "

$ws.Range("E35").Value = "  -1.29%  "

$ws.Range("D36").Value = "0.7073"
$ws.Range("E36").Value = "  +0.17%  "

$ws.Range("E37").Value = "  -3.22%  "

$ws.Range("D38").Value = "2.668"
$ws.Range("E38").Value = "  +1.02%  "

$ws.Range("D39").Value = "0.01852"
$ws.Range("E39").Value = "  -0.74%  "

$ws.Range("D40").Value = "2.729"
$ws.Range("E40").Value = "  +1.88%  "

$ws.Range("D41").Value = "1.145.71"
$ws.Range("E41").Value = "  +6.46%  "

$ws.Range("D42").Value = "0.9186"
$ws.Range("E42").Value = "  -1.87%  "

$ws.Range("D43").Value = "5.983"
$ws.Range("E43").Value = "  -0.26%  "

$ws.Range("D44").Value = "0.4287"

$ws.Range("D45").Value = "70.87"
$ws.Range("E45").Value = "  +0.94%  "

$ws.Range("D46").Value = "0.9994"
$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("D47").Value = "'102.70"
$ws.Range("E47").Value = "  +0.37%  "

$ws.Range("D48").Value = "1.773"
$ws.Range("E48").Value = "  +1.23%  "

$ws.Range("D49").Value = "2.006.00"
$ws.Range("E49").Value = "  +0.75%  "

$ws.Range("D50").Value = "9.182"
$ws.Range("E50").Value = "  +0.08%  "

$ws.Range("D51").Value = "6.977"
